$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (style) from H1 into the new I1/J1 header
# cells, then set their text values to the new header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-26: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 26; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
